# Update "want to go" counts (column F) that changed between scraper runs.
# Sheet "展览" (Exhibitions)
$wsExhibit = $excel.ActiveWorkbook.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 119
$wsExhibit.Range("F3").Value = 424

# Sheet "演出" (Performances)
$wsShow = $excel.ActiveWorkbook.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 68
$wsShow.Range("F3").Value = 27

# Sheet "全部类型" (All types, aggregated view of the sheets above)
$wsAll = $excel.ActiveWorkbook.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 119
$wsAll.Range("F3").Value = 68
$wsAll.Range("F4").Value = 424
$wsAll.Range("F8").Value = 27
